$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "Bloque"
$ws.Range("B1").Value = "Incidencia"
$ws.Range("C1").Value = "Fecha"
$ws.Range("D1").Value = "Hora"
$ws.Range("E1").Value = "Turno"
$ws.Range("F1").Value = "Hora de Reparación"
$ws.Range("G1").Value = "Tiempo de Reparación"
$ws.Range("H1").Value = "MTBF"

# --- Existing rows 2-4: fill in previously-missing F/G columns (already present, keep values) ---
$ws.Range("F2").Value = "12:18:22"
$ws.Range("G2").Value = "0:00:01"

$ws.Range("F3").Value = "12:18:27"
$ws.Range("G3").Value = "0:00:02"

$ws.Range("F4").Value = "12:18:39"
$ws.Range("G4").Value = "0:00:01"

# --- New rows 5-15 ---
$data = @(
  @("WC48 P5F", "AOI (fallo etiqueta)", "2024-05-29", "12:28:57", "Mañana", "12:28:58", "0:00:01", "N/A"),
  @("WC48 P5F", "Tornillo atascado en tolva", "2024-05-29", "12:29:15", "Mañana", "12:29:17", "0:00:02", "N/A"),
  @("WC48 P5F", "Power atascado en prensa, cuesta sacar", "2024-05-29", "12:29:27", "Mañana", "12:29:29", "0:00:02", "0.31 minutos"),
  @("WC48 P5F", "Cámara no detecta foam derecho", "2024-05-29", "12:29:45", "Mañana", "12:29:47", "0:00:02", "0.25 minutos"),
  @("WC48 P5F", "Robot no coge busbar", "2024-05-29", "12:30:06", "Mañana", "12:30:08", "0:00:02", "0.27 minutos"),
  @("WC48 P5F", "Cámara no detecta busbar", "2024-05-29", "12:30:10", "Mañana", "12:30:12", "0:00:02", "0.29 minutos"),
  @("WC49 P5H", "Tornillo atascado", "2024-05-29", "12:37:00", "Mañana", "12:37:03", "0:00:03", "N/A"),
  @("WC49 P5H", "Power atascado en prensa, cuesta sacar", "2024-05-29", "12:37:02", "Mañana", "12:37:04", "0:00:02", "N/A"),
  @("WC49 P5H", "La cámara no detecta Busbar", "2024-05-29", "12:37:21", "Mañana", "12:37:23", "0:00:02", "0.04 minutos"),
  @("WC49 P5H", "La cámara no detecta Busbar", "2024-05-29", "12:37:23", "Mañana", "12:37:25", "0:00:02", "0.17 minutos"),
  @("WC49 P5H", "La cámara no detecta Busbar", "2024-05-29", "12:37:35", "Mañana", "12:37:38", "0:00:03", "0.13 minutos")
)

$rowIndex = 5
foreach ($row in $data) {
  for ($col = 0; $col -lt $row.Length; $col++) {
    $cell = $ws.Cells.Item($rowIndex, $col + 1)
    $cell.NumberFormat = "@"
    $cell.Value = $row[$col]
  }
  $rowIndex++
}
